$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.892.74"
$ws.Range("E2").Value = "  -6.14%  "

# Row 3
$ws.Range("D3").Value = "2.987.49"
$ws.Range("E3").Value = "  -6.41%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.93"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "125.37"
$ws.Range("E6").Value = "  -8.64%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").Value = "2.984.92"
$ws.Range("E8").Value = "  -6.40%  "

# Row 9
$ws.Range("E9").Value = "  -2.53%  "

# Row 10
$ws.Range("E10").Value = "  -9.50%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.07"
$ws.Range("E11").Value = "  -5.53%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.440"
$ws.Range("E12").Value = "  -4.10%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000219"
$ws.Range("E13").Value = "  -9.41%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.49"
$ws.Range("E14").Value = "  -7.12%  "

# Row 15
$ws.Range("E15").Value = "  +0.82%  "

# Row 16
$ws.Range("D16").Value = "3.477.18"
$ws.Range("E16").Value = "  -6.25%  "

# Row 17
$ws.Range("D17").Value = "2.980.88"
$ws.Range("E17").Value = "  -6.61%  "

# Row 18
$ws.Range("D18").Value = "59.937.55"
$ws.Range("E18").Value = "  -6.03%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.45"
$ws.Range("E19").Value = "  -2.13%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "425.31"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.09"
$ws.Range("E21").Value = "  -6.44%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.667"
$ws.Range("E22").Value = "  -4.86%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.03"
$ws.Range("E23").Value = "  -8.60%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.94"
$ws.Range("E24").Value = "  -2.40%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.13"
$ws.Range("E25").Value = "  -4.85%  "

# Row 26
$ws.Range("E26").Value = "  +0.09%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.10%  "

# Row 28
$ws.Range("E28").Value = "  -6.33%  "

# Row 29
$ws.Range("E29").Value = "  -7.25%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.26"
$ws.Range("E30").Value = "  -7.26%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.19"
$ws.Range("E31").Value = "  -10.35%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.10"
$ws.Range("E32").Value = "  -9.13%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0938"
$ws.Range("E33").Value = "  -8.06%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.61"
$ws.Range("E34").Value = "  -4.99%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.929"
$ws.Range("E35").Value = "  -9.68%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.03"
$ws.Range("E36").Value = "  -3.28%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.06"
$ws.Range("E37").Value = "  -16.83%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0659"
$ws.Range("E38").Value = "  -10.91%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.34"
$ws.Range("E39").Value = "  +1.93%  "

# Row 40
$ws.Range("E40").Value = "  -10.58%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.108"

# Row 42
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "379.27"
$ws.Range("E42").Value = "  -4.99%  "

# Row 43
$ws.Range("D43").Value = "2.667.27"
$ws.Range("E43").Value = "  -4.93%  "

# Row 44
$ws.Range("E44").Value = "  -8.37%  "

# Row 46
$ws.Range("E46").Value = "  -7.70%  "

# Row 47
$ws.Range("E47").Value = "  -7.14%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.00"
$ws.Range("E48").Value = "  -7.12%  "

# Row 49
$ws.Range("E49").Value = "  -4.21%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.64"
$ws.Range("E50").Value = "  -7.81%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.00"
$ws.Range("E51").Value = "  -8.15%  "
